# Accounts.xlsx update: mark the visible rows of the first three test
# scenarios ("New Account", "New Case", "Edit Account") as Approved in the
# "Approved/Rejected" column (column I), matching the pattern already used
# for the later scenarios further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "Approved"
$ws.Range("I3").Value = "Approved"
$ws.Range("I4").Value = "Approved"
$ws.Range("I5").Value = "Approved"
$ws.Range("I6").Value = "Approved"
$ws.Range("I17").Value = "Approved"
$ws.Range("I25").Value = "Approved"

# Scroll the view back to the top of the sheet (column H) and select the
# newly-populated Approved/Rejected range.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 8
$ws.Range("I2:I61").Select()
